# Applies the TOYANOS.pptx edit:
#  1. Duplicates the last slide (old slide 9, "CONTACTO FACEBOOK, GITHUB.")
#     and moves the (still unedited) duplicate so it sits *before* the
#     original, becoming the new slide 9. The original slide (now slide 10,
#     still the same underlying part) is then retargeted with new content
#     ("Con el apoyo de los contribuyentes...").
#  2. Removes the "Con el apoyo de los contribuyentes..." bullet from slide 5
#     (Beneficios), since that text now lives on the new last slide.
#  3. Minor no-op "re-type" touch on slide 3's last bullet (cosmetic only).

$p = $ppt.ActivePresentation

# --- 1. Duplicate the last slide and reorder -----------------------------
$lastIndex = $p.Slides.Count            # 9
$original  = $p.Slides.Item($lastIndex)
$dup = $original.Duplicate()
$dup.MoveTo($lastIndex)                 # duplicate becomes slide 9, original slides down to 10

# The physical slide that keeps flowing content ("Con el apoyo...") is the
# one now at the end of the deck.
$target = $p.Slides.Item($p.Slides.Count)

$titleShape   = $target.Shapes.Item(1)
$contentShape = $target.Shapes.Item(2)

$titleShape.Name   = "Título 1"
$contentShape.Name = "Marcador de contenido 2"

$titleShape.TextFrame.TextRange.Text = "Con el apoyo de los contribuyentes"
$titleShape.TextFrame.TextRange.Paragraphs(1,1).LanguageID = "es-MX"

$contentShape.TextFrame.TextRange.Text = "Con el apoyo de los contribuyentes, ampliaremos la base de datos, su numero de`rusuarios, y de igual manera la ligaremos a instituciones confiables.`r"

$contentTr = $contentShape.TextFrame.TextRange
for ($i = 1; $i -le $contentTr.Paragraphs().Count; $i++) {
    $para = $contentTr.Paragraphs($i, 1)
    $para.LanguageID = "es-MX"
    $para.Font.Bold = $false
}

# --- 2. Slide 5 ("Beneficios"): drop the paragraph that moved to slide 10 -
$s5 = $p.Slides.Item(5)
$s5content = $s5.Shapes.Item(2)
$s5tr = $s5content.TextFrame.TextRange
$s5tr.Paragraphs(2, 1).Delete()

# --- 3. Slide 3 ("Proyecto Troyanos"): cosmetic re-touch of last bullet --
$s3 = $p.Slides.Item(3)
$s3content = $s3.Shapes.Item(2)
$s3tr = $s3content.TextFrame.TextRange
$lastPara = $s3tr.Paragraphs($s3tr.Paragraphs().Count, 1)
$lastPara.Text = "Es fácil de usar y es agradable a la vista"
